$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-22 20:38:52"
$wsOverview.Columns.Item(5).ColumnWidth = 16.4
$wsOverview.Columns.Item(6).ColumnWidth = 16.4

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-22 20:38:47"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.4

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-22 20:38:52"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.4
